$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "https://en.wikipedia.org/wiki/Main_Page"
$ws.Range("A5").Value = "testing 15 / 15"
